# Append the 02/03/2026 daily snapshot as row 71 of Sheet1.
# (mirrors: "Update profit files after running on 2026-02-03")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a plain date-formatted label (literal text, not a real
# Excel date serial) - every existing row in this sheet stores it that way.
# Leading "'" forces text entry so "02/03/2026" isn't auto-converted to a
# date serial; resetting the style back to Normal afterwards drops the
# quote-prefix formatting that the text-forced entry would otherwise leave
# behind, matching the unstyled cells used by every other row.
$ws.Range("A71").Value = "'02/03/2026"
$ws.Range("A71").Style = "Normal"

$ws.Range("B71").Value = 9956.700000000001
$ws.Range("C71").Value = 0.2526948899720284
$ws.Range("D71").Value = 0.7473051100279716
$ws.Range("E71").Value = -296.87
$ws.Range("F71").Value = -37.82
$ws.Range("G71").Value = -23401.06
$ws.Range("H71").Value = -75.87
$ws.Range("I71").Value = -786.84
$ws.Range("J71").Value = -23.82
$ws.Range("K71").Value = -24187.9
$ws.Range("L71").Value = -70.84
